# Apply updated odds values per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("O2").Value = 1.37

# Row 5
$ws.Range("F5").Value = 3.7
$ws.Range("G5").Value = 4.6
$ws.Range("H5").Value = 1.85
$ws.Range("I5").Value = 2.22
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.3
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.25
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 1.88
$ws.Range("Q5").Value = 1.87
$ws.Range("R5").Value = 1.33
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 1.75
$ws.Range("U5").Value = 2.04
$ws.Range("V5").Value = 1.82
$ws.Range("W5").Value = 1.27
$ws.Range("X5").Value = 17.5
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 14.5
$ws.Range("AA5").Value = 27
$ws.Range("AB5").Value = 18
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 11.5
$ws.Range("AE5").Value = 27
$ws.Range("AF5").Value = 38
$ws.Range("AG5").Value = 20
$ws.Range("AH5").Value = 22
$ws.Range("AI5").Value = 44
$ws.Range("AK5").Value = 65
$ws.Range("AN5").Value = 65
$ws.Range("AO5").Value = 16

# Row 7
$ws.Range("F7").Value = 6.2
$ws.Range("G7").Value = 9.199999999999999
$ws.Range("H7").Value = 1.44
$ws.Range("I7").Value = 1.54
$ws.Range("J7").Value = 4.4
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 1.27
$ws.Range("M7").Value = 1.05
$ws.Range("N7").Value = 3.45
$ws.Range("O7").Value = 1.26
$ws.Range("P7").Value = 1.9
$ws.Range("Q7").Value = 1.71
$ws.Range("R7").Value = 1.37
$ws.Range("S7").Value = 2.86
$ws.Range("T7").Value = 2.02
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 2.86
$ws.Range("AF7").Value = 80

# Row 9
$ws.Range("J9").Value = 3.5
$ws.Range("V9").Value = 1.38

# Row 10
$ws.Range("U10").Value = 2
$ws.Range("AE10").Value = 17.5
$ws.Range("AN10").Value = 95

# Row 11
$ws.Range("G11").Value = 13
$ws.Range("I11").Value = 1.3

# Row 12
$ws.Range("G12").Value = 6.4
$ws.Range("H12").Value = 1.67
$ws.Range("I12").Value = 1.88
$ws.Range("J12").Value = 2.98
$ws.Range("O12").Value = 1.2
$ws.Range("Q12").Value = 1.61
$ws.Range("V12").Value = 2.12
$ws.Range("W12").Value = 1.18

# Row 13
$ws.Range("G13").Value = 1.35
$ws.Range("J13").Value = 4.5
$ws.Range("V13").Value = 1.06
